$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Legend column (D) currently holds, top to bottom starting at row 3:
#   D3 red    "Solicitada"
#   D4 blue   "Aprovada"
#   D5 green  "Executada"
#   D6..D14   empty bordered cells (D9 has the special "no bottom border" style)
#
# We need to insert a new first legend entry "Planejada" (yellow) above the
# existing three, which pushes them down by one row, and append one more
# blank bordered row at the bottom (row 15). Shift bottom-up so sources
# aren't clobbered before they're read, using Range.Copy(Destination) which
# copies both value and formatting in one shot.

$ws.Range("D5").Copy($ws.Range("D6"))
$ws.Range("D4").Copy($ws.Range("D5"))
$ws.Range("D3").Copy($ws.Range("D4"))

# D3 becomes the new legend entry: same bordered style as the others, but a
# new yellow fill, and the new text.
$ws.Range("D3").Value = "Planejada"
$ws.Range("D3").Interior.Color = 65535

# The blank-cell block below the legend used to end its "no bottom border"
# marker row at D9; after the insertion above it that marker moves to D10.
$ws.Range("D9").Copy($ws.Range("D10"))
$ws.Range("D8").Copy($ws.Range("D9"))

# Append the new trailing blank bordered row 15 (column D only).
$ws.Range("D14").Copy($ws.Range("D15"))

# Selection moved from E5 to D5 in the saved view state.
$ws.Range("D5").Select()
